# The legacy GSC export "Chart" sheet carried a placeholder row (row 2)
# for 2025-10-12 with no data yet (empty strings in columns B/C).
# That date has since been backfilled upstream, so the stale placeholder
# row is removed here; Excel shifts every subsequent row up by one and
# the sheet's used range shrinks by a row (A1:D83 -> A1:D82).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows(2).Delete()
